$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert two new blank columns at K:L (shifts old K..P to M..R)
$ws.Columns("K:L").Insert()

# Step 2: Update header row (row 1) for new/moved columns
$ws.Range("K1").Value = "etrucks"
$ws.Range("L1").Value = "ebestel"
$ws.Range("Q1").Value = "jaarkilometrage"
$ws.Range("R1").Value = "jaarkilometrage_truck"
$ws.Range("S1").Value = "jaarkilometrage_bestel"

# Make sure the newly appended header cells (Q1:S1) carry the same header
# formatting (bold, centered, thin border) as the rest of row 1, since they
# were appended beyond the sheet's original used range and don't
# automatically inherit that formatting.
$ws.Range("P1").Copy() | Out-Null
$ws.Range("Q1:S1").PasteSpecial(-4122) | Out-Null

# Step 3: Fix E/F (trucks/bestelbussen) values that changed for a few rows
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("E5").Value = 38
$ws.Range("F5").Value = 5
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 3
$ws.Range("E30").Value = 6
$ws.Range("F30").Value = 4

# Step 4: Populate new K:L columns (etrucks, ebestel) for rows 2-33
$klArr = New-Object 'object[,]' 32,2
$klArr[0,0] = 70
$klArr[0,1] = 0
$klArr[1,0] = 0
$klArr[1,1] = 0
$klArr[2,0] = 0
$klArr[2,1] = 0
$klArr[3,0] = 2
$klArr[3,1] = 1
$klArr[4,0] = 2
$klArr[4,1] = 4
$klArr[5,0] = 0
$klArr[5,1] = 0
$klArr[6,0] = 0
$klArr[6,1] = 0
$klArr[7,0] = 0
$klArr[7,1] = 0
$klArr[8,0] = 0
$klArr[8,1] = 0
$klArr[9,0] = 0
$klArr[9,1] = 0
$klArr[10,0] = 0
$klArr[10,1] = 0
$klArr[11,0] = 0
$klArr[11,1] = 0
$klArr[12,0] = 0
$klArr[12,1] = 0
$klArr[13,0] = 0
$klArr[13,1] = 0
$klArr[14,0] = 0
$klArr[14,1] = 0
$klArr[15,0] = 0
$klArr[15,1] = 0
$klArr[16,0] = 0
$klArr[16,1] = 0
$klArr[17,0] = 0
$klArr[17,1] = 0
$klArr[18,0] = 0
$klArr[18,1] = 0
$klArr[19,0] = 0
$klArr[19,1] = 0
$klArr[20,0] = 0
$klArr[20,1] = 0
$klArr[21,0] = 0
$klArr[21,1] = 0
$klArr[22,0] = 0
$klArr[22,1] = 0
$klArr[23,0] = 0
$klArr[23,1] = 0
$klArr[24,0] = 0
$klArr[24,1] = 0
$klArr[25,0] = 0
$klArr[25,1] = 0
$klArr[26,0] = 0
$klArr[26,1] = 0
$klArr[27,0] = 0
$klArr[27,1] = 0
$klArr[28,0] = 3
$klArr[28,1] = 0
$klArr[29,0] = 0
$klArr[29,1] = 0
$klArr[30,0] = 0
$klArr[30,1] = 0
$klArr[31,0] = 0
$klArr[31,1] = 0
$ws.Range("K2:L33").Value = $klArr

# Step 5: Populate new Q:S columns (jaarkilometrage, _truck, _bestel) for rows 2-33
$qrsArr = New-Object 'object[,]' 32,3
$qrsArr[0,0] = 18000
$qrsArr[0,1] = 54505
$qrsArr[0,2] = 18000
$qrsArr[1,0] = 18000
$qrsArr[1,1] = 54505
$qrsArr[1,2] = 18000
$qrsArr[2,0] = 20000
$qrsArr[2,1] = 54505
$qrsArr[2,2] = 20000
$qrsArr[3,0] = 22500
$qrsArr[3,1] = 22500
$qrsArr[3,2] = 18000
$qrsArr[4,0] = 18000
$qrsArr[4,1] = 54505
$qrsArr[4,2] = 18000
$qrsArr[5,0] = 18000
$qrsArr[5,1] = 54505
$qrsArr[5,2] = 18000
$qrsArr[6,0] = 18000
$qrsArr[6,1] = 54505
$qrsArr[6,2] = 18000
$qrsArr[7,0] = 18000
$qrsArr[7,1] = 54505
$qrsArr[7,2] = 18000
$qrsArr[8,0] = 30000
$qrsArr[8,1] = 54505
$qrsArr[8,2] = 30000
$qrsArr[9,0] = 30000
$qrsArr[9,1] = 54505
$qrsArr[9,2] = 30000
$qrsArr[10,0] = 30000
$qrsArr[10,1] = 54505
$qrsArr[10,2] = 30000
$qrsArr[11,0] = 15000
$qrsArr[11,1] = 54505
$qrsArr[11,2] = 15000
$qrsArr[12,0] = 30000
$qrsArr[12,1] = 54505
$qrsArr[12,2] = 30000
$qrsArr[13,0] = 18000
$qrsArr[13,1] = 54505
$qrsArr[13,2] = 18000
$qrsArr[14,0] = 18000
$qrsArr[14,1] = 54505
$qrsArr[14,2] = 18000
$qrsArr[15,0] = 18000
$qrsArr[15,1] = 54505
$qrsArr[15,2] = 18000
$qrsArr[16,0] = 18000
$qrsArr[16,1] = 54505
$qrsArr[16,2] = 18000
$qrsArr[17,0] = 30000
$qrsArr[17,1] = 54505
$qrsArr[17,2] = 30000
$qrsArr[18,0] = 18000
$qrsArr[18,1] = 54505
$qrsArr[18,2] = 18000
$qrsArr[19,0] = 18000
$qrsArr[19,1] = 54505
$qrsArr[19,2] = 18000
$qrsArr[20,0] = 30000
$qrsArr[20,1] = 54505
$qrsArr[20,2] = 30000
$qrsArr[21,0] = 137500
$qrsArr[21,1] = 137500
$qrsArr[21,2] = 30000
$qrsArr[22,0] = 18000
$qrsArr[22,1] = 54505
$qrsArr[22,2] = 18000
$qrsArr[23,0] = 18000
$qrsArr[23,1] = 54505
$qrsArr[23,2] = 18000
$qrsArr[24,0] = 30000
$qrsArr[24,1] = 54505
$qrsArr[24,2] = 30000
$qrsArr[25,0] = 18000
$qrsArr[25,1] = 54505
$qrsArr[25,2] = 18000
$qrsArr[26,0] = 20000
$qrsArr[26,1] = 54505
$qrsArr[26,2] = 20000
$qrsArr[27,0] = 18000
$qrsArr[27,1] = 54505
$qrsArr[27,2] = 18000
$qrsArr[28,0] = 18000
$qrsArr[28,1] = 54505
$qrsArr[28,2] = 18000
$qrsArr[29,0] = 18000
$qrsArr[29,1] = 54505
$qrsArr[29,2] = 18000
$qrsArr[30,0] = 18000
$qrsArr[30,1] = 54505
$qrsArr[30,2] = 18000
$qrsArr[31,0] = 18000
$qrsArr[31,1] = 54505
$qrsArr[31,2] = 18000
$ws.Range("Q2:S33").Value = $qrsArr

Write-Output "Done"